# Generate Report for Handback
# Updates the handoff/handback timestamps produced for the 853a25db... file
# (first data row on each sheet) as part of a fresh report-generation run.
# The ac1ba60d... file's row (second data row) is left untouched.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview: "Latest HO Xliff Generate Date" for 853a25db row
$overview.Range("G2").Value = "2016-09-03 19:01:31"

# zh-cn: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for 853a25db row
$zhcn.Range("H2").Value = "2016-09-03 19:01:27"
$zhcn.Range("K2").Value = "2016-09-03 19:01:45"

# de-de: "Correspond Handback DateTime" for 853a25db row
$dede.Range("K2").Value = "2016-09-03 19:01:52"
